$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Area closures" cell (B8) now carries the same text as California/Oregon's
# "Area closure, evisceration order (area closure only in 2015)", so the shared
# string "Area closures" disappears from the table and B8 adopts the wrap-text
# style already used by C8/D8.
$ws.Range("B8").Value = "Area closure, evisceration order`n(area closure only in 2015)"
$ws.Range("B8").WrapText = $true
$ws.Range("B8").VerticalAlignment = -4160  # xlTop

# Column B is widened (and loses its "best fit" autosize flag) to comfortably
# show the new, longer text.
$ws.Columns.Item(2).ColumnWidth = 27.5

# Move the saved selection.
$ws.Range("B10").Select()
